$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.911.85'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.093.11'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '543.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.28%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.088.38'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('E9').Value = '  +1.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.54'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.77%  '
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.457'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('E13').Value = '  +4.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.79'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.589.11'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.970.10'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.56%  '
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.091.30'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.66'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '479.38'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.11%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.700'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.32'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.06'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.73%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('E31').Value = '  -2.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.15'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '57.20'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.84%  '
$ws.Range('E34').Value = '  -5.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '494.01'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.37'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.02'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.250.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.57%  '
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('E40').Value = '  +0.75%  '
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.68'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.12'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '124.29'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '25.13'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.55%  '
$ws.Range('E48').Value = '  -1.67%  '
$ws.Range('E49').Value = '  +7.24%  '
$ws.Range('E50').Value = '  +1.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.40'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.00%  '
